# Auto-generated: update cryptos list values (prices, volumes, and re-ranked rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "56.497.31"
$ws.Range("E2").Value = "  +10.61%  "

# Row 3: Ethereum
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "3.245.80"
$ws.Range("E3").Value = "  +6.18%  "

# Row 4: TetherUSD
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5: BNB
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'398.40"
$ws.Range("E5").Value = "  +3.02%  "

# Row 6: Solana
$ws.Range("B6").Value = "Solana"
$ws.Range("C6").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D6").Value = "'111.34"
$ws.Range("E6").Value = "  +8.88%  "

# Row 7: XRP
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.560"
$ws.Range("E7").Value = "  +4.60%  "

# Row 8: USDC
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.12%  "

# Row 9: Cardano
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +7.65%  "

# Row 10: Avalanche
$ws.Range("B10").Value = "Avalanche"
$ws.Range("C10").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D10").Value = "'39.38"
$ws.Range("E10").Value = "  +7.55%  "

# Row 11: Dogecoin
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0947"
$ws.Range("E11").Value = "  +11.67%  "

# Row 12: TRON
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.141"
$ws.Range("E12").Value = "  +2.34%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.755.47"
$ws.Range("E13").Value = "  +5.66%  "

# Row 14: Polkadot
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'8.13"
$ws.Range("E14").Value = "  +6.00%  "

# Row 15: Chainlink
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'19.18"
$ws.Range("E15").Value = "  +4.89%  "

# Row 16: WrappedEther
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.252.61"
$ws.Range("E16").Value = "  +6.14%  "

# Row 17: Polygon
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'1.05"
$ws.Range("E17").Value = "  +6.41%  "

# Row 18: Uniswap
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'11.05"
$ws.Range("E18").Value = "  +3.97%  "

# Row 19: WrappedBTC
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "56.314.89"
$ws.Range("E19").Value = "  +10.14%  "

# Row 20: ImmutableX
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").Value = "'3.35"
$ws.Range("E20").Value = "  +4.75%  "

# Row 21: ShibaInu
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000104"
$ws.Range("E21").Value = "  +9.31%  "

# Row 22: InternetComputer(DFINITY)
$ws.Range("B22").Value = "InternetComputer(DFINITY)"
$ws.Range("C22").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D22").Value = "'13.17"
$ws.Range("E22").Value = "  +7.62%  "

# Row 23: BitcoinCash
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'299.25"
$ws.Range("E23").Value = "  +13.36%  "

# Row 24: Litecoin
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'75.80"
$ws.Range("E24").Value = "  +8.79%  "

# Row 25: PancakeSwap
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'3.24"
$ws.Range("E25").Value = "  +3.63%  "

# Row 26: Filecoin
$ws.Range("B26").Value = "Filecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D26").Value = "'8.22"
$ws.Range("E26").Value = "  +4.39%  "

# Row 27: EthereumClassic
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'28.20"
$ws.Range("E27").Value = "  +4.64%  "

# Row 28: LEO
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "'4.34"
$ws.Range("E28").Value = "  +3.99%  "

# Row 29: RenderToken
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.36"
$ws.Range("E29").Value = "  +2.29%  "

# Row 30: Kaspa
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.171"
$ws.Range("E30").Value = "  +4.74%  "

# Row 31: Dai
$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.06%  "

# Row 32: Hedera
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.111"
$ws.Range("E32").Value = "  +6.03%  "

# Row 33: Cosmos
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.17"
$ws.Range("E33").Value = "  +7.53%  "

# Row 34: InjectiveProtocol
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'36.94"
$ws.Range("E34").Value = "  +4.61%  "

# Row 35: VeChain
$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "'0.0491"
$ws.Range("E35").Value = "  +4.49%  "

# Row 36: Toncoin
$ws.Range("B36").Value = "Toncoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D36").Value = "'2.13"
$ws.Range("E36").Value = "  +3.12%  "

# Row 37: OKB
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'51.41"
$ws.Range("E37").Value = "  +2.79%  "

# Row 38: LidoDAOToken
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "'3.55"
$ws.Range("E38").Value = "  +5.19%  "

# Row 39: FirstDigitalUSD
$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40: Stacks
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "'3.09"
$ws.Range("E40").Value = "  +25.64%  "

# Row 41: Monero
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "'135.08"
$ws.Range("E41").Value = "  +4.35%  "

# Row 42: Celestia
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'17.62"
$ws.Range("E42").Value = "  +7.78%  "

# Row 43: ARBITRUM
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "'1.93"
$ws.Range("E43").Value = "  +5.49%  "

# Row 44: NEARProtocol
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'4.00"
$ws.Range("E44").Value = "  +5.36%  "

# Row 45: Stellar
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.120"
$ws.Range("E45").Value = "  +4.63%  "

# Row 46: TheGraph
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.284"
$ws.Range("E46").Value = "  -1.19%  "

# Row 47: EnergySwap
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'22.31"
$ws.Range("E47").Value = "  +3.15%  "

# Row 48: ThetaToken
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").Value = "'2.19"
$ws.Range("E48").Value = "  +56.68%  "

# Row 49: Maker
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.144.10"
$ws.Range("E49").Value = "  +4.18%  "

# Row 50: WEMIXToken
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'2.09"
$ws.Range("E50").Value = "  +0.73%  "

# Row 51: ApeXProtocol
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'2.44"
$ws.Range("E51").Value = "  -1.62%  "

